# edit.ps1 - applies the "midway report.docx" revision:
#   1) Merge the three runs (with the now-removed __DdeLink__ bookmark)
#      that spell out "log((total # of views for uploader's videos ...)"
#      back into a single run of plain text.
#   2) Append a new "Next Steps" section (heading, body paragraph, and a
#      trailing empty bold paragraph) after the existing "DLKfsdf" line.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Collapse the "log((total # of views for uploader's videos ...))"
#    run/bookmark/run/bookmark/run sequence into one run.
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(
    "log((total # of views for ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Extend the found range so it spans all three original runs
    # (through the end of "...minus 1))"), swallowing the bookmark
    # start/end markers that sat between them.
    $r.MoveEnd(1, 92)

    # Assigning identical text is a no-op in this engine, so first
    # stamp a placeholder, then set the real text - this guarantees
    # the whole span is rebuilt as a single run (and the now-orphaned
    # bookmark tags are dropped) while keeping the original run
    # formatting (b=false, bCs=false, u=none).
    $r.Text = "__PLACEHOLDER__"
    $r.Text = "log((total # of views for uploader's videos " + `
        [char]0x2013 + " views for this video) / (# of videos uploaded by the uploader, minus 1))"
}

# ---------------------------------------------------------------------
# 2) Append the "Next Steps" section after the "DLKfsdf" paragraph.
# ---------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$headingXml = '<w:p ' + $wNs + '>' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="Normal"/>' + `
        '<w:spacing w:before="0" w:after="144"/>' + `
        '<w:jc w:val="left"/>' + `
        '<w:rPr><w:b/><w:bCs/><w:u w:val="none"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr><w:b/><w:bCs/><w:u w:val="none"/></w:rPr>' + `
        '<w:t>Next Steps</w:t>' + `
    '</w:r>' + `
'</w:p>'

$bodyXml = '<w:p ' + $wNs + '>' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="Normal"/>' + `
        '<w:spacing w:before="0" w:after="144"/>' + `
        '<w:jc w:val="left"/>' + `
        '<w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:u w:val="none"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:u w:val="none"/></w:rPr>' + `
        '<w:tab/>' + `
        '<w:t>Based on the results observed, we need to LKDJFLSJFLKSD</w:t>' + `
    '</w:r>' + `
'</w:p>'

$trailerXml = '<w:p ' + $wNs + '>' + `
    '<w:pPr>' + `
        '<w:pStyle w:val="Normal"/>' + `
        '<w:spacing w:before="0" w:after="144"/>' + `
        '<w:jc w:val="left"/>' + `
        '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r>' + `
        '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
    '</w:r>' + `
'</w:p>'

$p1 = $d.Paragraphs.Add()
$p1.Range.InsertXML($headingXml)

$p2 = $d.Paragraphs.Add()
$p2.Range.InsertXML($bodyXml)

$p3 = $d.Paragraphs.Add()
$p3.Range.InsertXML($trailerXml)

Write-Output "Done. Paragraphs.Count=$($d.Paragraphs.Count)"
